$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selenium")

$ws.Range("A1").Value = "App"
$ws.Range("B1").Value = "TestCase ID"
$ws.Range("C1").Value = "Test Case Name"
$ws.Range("D1").Value = "Test Data"
$ws.Range("E1").Value = "Pre Conditions"
$ws.Range("F1").Value = "Steps to Perform"
$ws.Range("G1").Value = "Expected Results"

$ws.Range("A2").Value = "Green Kart - Shopping Page"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Verify the company logo is displayed in header when navigating to page."
$ws.Range("D2").Value = "N/A"
$ws.Range("I2").Value = "Websites"
$ws.Range("J2").Value = "https://rahulshettyacademy.com/seleniumPractise/#/"

$ws.Range("A3").Value = "Green Kart - Shopping Page"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Verify product list shows the expected products when using a full product search."
$ws.Range("D3").Value = "Product to search for."
$ws.Range("J3").Value = "https://rahulshettyacademy.com/AutomationPractice/"

$ws.Range("A4").Value = "Green Kart - Shopping Page"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Verify product list shows products when using a parital product search."
$ws.Range("D4").Value = "Partial string to search for."
$ws.Range("J4").Value = "https://rahulshettyacademy.com/dropdownsPractise/"

$ws.Range("A5").Value = "Green Kart - Shopping Page"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Verify product list shows the expected error when using an invalid product search."
$ws.Range("D5").Value = "Invalid string to search for."
$ws.Range("J5").Value = "https://rahulshettyacademy.com/angularpractice/"

$ws.Range("A6").Value = "Green Kart - Shopping Page"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Verify product list returns all products after search has been cleared."
$ws.Range("D6").Value = "Product list"

$ws.Range("A7").Value = "Green Kart - Shopping Page"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = "Verify product list shows the correct price for each item."
$ws.Range("D7").Value = "Product list with prices"

$ws.Range("A8").Value = "Green Kart - Shopping Page"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "Verify product list shows products in the correct order."

$ws.Range("A9").Value = "Green Kart - Shopping Page"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = "Verify Top Deals link navigates to the correct page."
$ws.Range("D9").Value = "N/A"

$ws.Range("A10").Value = "Green Kart - Shopping Page"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = "Verify Flight Booking link navigates to the correct page."
$ws.Range("D10").Value = "N/A"

$ws.Range("A11").Value = "Green Kart - Shopping Page"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "Verify limited offer link navigates to the correct page."

$ws.Range("A12").Value = "Green Kart - Shopping Page"
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = "Verify Proceed to Checkout button in the cart is disabled when cart is empty."

$ws.Range("A13").Value = "Green Kart - Shopping Page"
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = "Verify cart shows the cart is empty message when cart is empty."

$ws.Range("A14").Value = "Green Kart - Shopping Page"
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = "Verify product shows the added message when it is added to the cart."

$ws.Range("A15").Value = "Green Kart - Shopping Page"
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = "Verify cart shows product when a product has been added to the cart."

$ws.Range("A16").Value = "Green Kart - Shopping Page"
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = "Verify cart shows the correct quantity if multiple of one item are added to the cart."

$ws.Range("A17").Value = "Green Kart - Shopping Page"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = "Verify cart shows all items in cart when multiple different items have been added."

$ws.Range("A18").Value = "Green Kart - Shopping Page"
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = "Verify cart still shows added items after refreshing the page."

$ws.Range("A19").Value = "Green Kart - Shopping Page"
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = "Verify cart no longer shows items that have been deleted from the cart."

$ws.Range("A20").Value = "Green Kart - Shopping Page"
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = "Verify cart shows the correct total for items that have been added to the cart."

$ws.Range("A21").Value = "Green Kart - Checkout Page"
$ws.Range("B21").Value = 19
$ws.Range("C21").Value = "Verify table lists all products that were added to the cart."

$ws.Range("A22").Value = "Green Kart - Checkout Page"
$ws.Range("B22").Value = 20
$ws.Range("C22").Value = "Verify table has the correct quantity for each product that was added to the cart."

$ws.Range("A23").Value = "Green Kart - Checkout Page"
$ws.Range("B23").Value = 21
$ws.Range("C23").Value = "Verify table has the correct price for each product that was added to the cart."

$ws.Range("A24").Value = "Green Kart - Checkout Page"
$ws.Range("B24").Value = 22
$ws.Range("C24").Value = "Veirfy table is calculating the total correctly for each product added to the cart."

$ws.Range("A25").Value = "Green Kart - Checkout Page"
$ws.Range("B25").Value = 23
$ws.Range("C25").Value = "Verify total is displayed and correct for all items in the cart."

$ws.Range("A26").Value = "Green Kart - Checkout Page"
$ws.Range("B26").Value = 24
$ws.Range("C26").Value = "Verify coupon code shows correct error message when an invalid code is entered."

$ws.Range("A27").Value = "Green Kart - Checkout Page"
$ws.Range("B27").Value = 25
$ws.Range("C27").Value = "Verify coupon code shows correct success message when a valid code is entered."

$ws.Range("A28").Value = "Green Kart - Checkout Page"
$ws.Range("B28").Value = 26
$ws.Range("C28").Value = "Verify total after discount is calculated and displayed correctly when entering a valid coupon."

$ws.Range("A29").Value = "Green Kart - Checkout Page"
$ws.Range("B29").Value = 27
$ws.Range("C29").Value = "Verify place order button takes user to the shipping page."

$ws.Range("A30").Value = "Green Kart - Shipping Page"
$ws.Range("B30").Value = 28
$ws.Range("C30").Value = "Verify country dropdown allows a country to be selected."

$ws.Range("A31").Value = "Green Kart - Shipping Page"
$ws.Range("B31").Value = 29
$ws.Range("C31").Value = "Verify error message is displayed if terms and conditions are not agreed to."

$ws.Range("A32").Value = "Green Kart - Shipping Page"
$ws.Range("B32").Value = 30
$ws.Range("C32").Value = "Verify success message is displayed upon successfully placing order."

$ws.Range("A33").Value = "Green Kart - Shipping Page"
$ws.Range("B33").Value = 31
$ws.Range("C33").Value = "Verify shopping page is shown after completing a successful order."

$ws.Columns.Item(1).ColumnWidth = 24.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 90.66666666666667
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668
$ws.Columns.Item(5).ColumnWidth = 26.5
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668
$ws.Columns.Item(7).ColumnWidth = 29.166666666666668

$ws.Range("C33").Select()
